$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet updates
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value2 = -0.16      # Total P&L %
$summary.Range("B6").Value2 = 26         # Total Trades
$summary.Range("B9").Value2 = 46.15      # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet updates (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value2 = 26          # Trades
$status.Range("G5").Value2 = 46.15       # Win Rate %

# ---------------------------------------------------------------------------
# 3. Helper that appends trade #26 (row 27) to a trade-log sheet
# ---------------------------------------------------------------------------
function Add-Trade26Row {
    param($ws)

    $ws.Range("A27").Value2 = 26

    # Date/Time columns look like numbers/dates to Excel's parser, so force
    # them to be stored as text before assigning, just like the existing
    # rows above them.
    $ws.Range("B27").NumberFormat = "@"
    $ws.Range("B27").Value2 = "2026-02-17"

    $ws.Range("C27").Value2 = "20:07:31"
    $ws.Range("D27").Value2 = "MarketMaking"
    $ws.Range("E27").Value2 = "UP"
    $ws.Range("F27").Value2 = 0.97
    $ws.Range("G27").Value2 = 0.97
    $ws.Range("H27").Value2 = "CLOSED"
    $ws.Range("I27").Value2 = 0
    $ws.Range("J27").Value2 = 0
    $ws.Range("K27").Value2 = 99.8
    $ws.Range("L27").Value2 = 0
    $ws.Range("M27").Value2 = 0
    $ws.Range("N27").Value2 = 0.6
    $ws.Range("O27").Value2 = "Normal spread capture: 19600 bps"
    $ws.Range("P27").Value2 = "early_exit"
    $ws.Range("Q27").Value2 = 0.13
}

# ---------------------------------------------------------------------------
# 4. Apply the new trade row to both "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade26Row $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade26Row $marketMaking
